$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = "5056780-46.2019.8.21.0001"
$ws.Range("B2").Value = "0140462-81.2019.8.21.0001"

# Update row 3
$ws.Range("A3").Value = "5000559-78.2019.8.21.0054"
$ws.Range("B3").Value = "0003337-09.2019.8.21.0054"

# Update row 4
$ws.Range("A4").Value = "5009986-45.2011.8.21.0001"
$ws.Range("B4").Value = "0323125-76.2011.8.21.0001"

# Delete rows 5 through 13 (entire rows), which removes the remaining data
$ws.Range("A5:C13").EntireRow.Delete() | Out-Null
